$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")
$ws.Rows.Item(21).Delete()
